$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 14 is a duplicate of the current (pre-edit) row 13 contents,
# since row 13 itself is about to be updated with new values below.
$ws.Range("A14").Value2 = $ws.Range("A13").Value2
$ws.Range("B14").Value2 = $ws.Range("B13").Value2
$ws.Range("C14").Value2 = $ws.Range("C13").Value2
$ws.Range("D14").Value2 = $ws.Range("D13").Value2
$ws.Range("E14").Value2 = $ws.Range("E13").Value2
$ws.Range("F14").Value2 = $ws.Range("F13").Value2
$ws.Range("G14").Value2 = $ws.Range("G13").Value2
$ws.Range("H14").Value2 = $ws.Range("H13").Value2
$ws.Range("I14").Value2 = $ws.Range("I13").Value2
$ws.Range("J14").Value2 = $ws.Range("J13").Value2
$ws.Range("K14").Value2 = $ws.Range("K13").Value2
$ws.Range("L14").Value2 = $ws.Range("L13").Value2
$ws.Range("M14").Value2 = $ws.Range("M13").Value2
$ws.Range("N14").Value2 = $ws.Range("N13").Value2
$ws.Range("O14").Value2 = $ws.Range("O13").Value2
$ws.Range("P14").Value2 = $ws.Range("P13").Value2
$ws.Range("Q14").Value2 = $ws.Range("Q13").Value2
$ws.Range("R14").Value2 = $ws.Range("R13").Value2
$ws.Range("S14").Value2 = $ws.Range("S13").Value2
$ws.Range("T14").Value2 = $ws.Range("T13").Value2

# Copy D13's number format (date style) onto D14 without touching its value.
$ws.Range("D13").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Now update row 13 with its new values.
$ws.Range("D13").Value2 = 44627
$ws.Range("M13").Value2 = 45
